$d = $word.ActiveDocument

# Useful characters
$ldq = [char]0x201C
$rdq = [char]0x201D
$apos = [char]0x2019
$cr = [char]13

# ---------------------------------------------------------------------------
# Step 1: extend the "We used the framework..." sentence and append a new
# sentence (second run) to the same paragraph, preserving the existing
# i=false/iCs=false/lang=en-US run formatting.
# ---------------------------------------------------------------------------
$old1 = $ldq + "layers" + $rdq + " to understand and store our intermediate results. "
$new1 = $ldq + "layers" + $rdq + " and " + $ldq + "subproblems" + $rdq + " to understand and store our intermediate results. We define a subproblem as a collection of a " + $ldq + "cost" + $rdq + " or length of a path, and the path itself. Layers are used to group series of subproblems together in a way to make the intermediate answers accessible for later calculations.  "
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

$p19 = $d.Paragraphs.Item(19)
$r19 = $p19.Range
$insStart = $r19.End - 1
$insPoint = $d.Range($insStart, $insStart)
$run2Text = "Entries in a layer are accessed by providing a layer index as well as a tuple describing the index of the city being considered and a set of visited cities."
$insPoint.InsertAfter($run2Text)

# Force the new text into its own run (splitting from run 1) while keeping
# the inherited i/iCs/lang triple intact.
$run2Range = $d.Range($insStart, $insStart + $run2Text.Length)
$run2Range.Font.Italic = $true
$run2Range.Font.Italic = $false

Write-Host ("Paragraph 19 now: " + $d.Paragraphs.Item(19).Range.Text)

# ---------------------------------------------------------------------------
# Step 2: insert the new paragraphs. We anchor on the paragraph that is
# currently empty (directly before "Acknowledgements") since text inserted
# at its start inherits a *clean* (no direct formatting) run, which matches
# most of the new runs we need to create.
# ---------------------------------------------------------------------------
$anchor = $d.Paragraphs.Item(20)
$anchorIns = $d.Range($anchor.Range.Start, $anchor.Range.Start)

$h1Text  = "2.2.1 The 0th layer"
$bodyA1  = "We populated our 0th layer by selecting an arbitrary starting city and storing subproblems relating the cost and route of all the connections coming from that city.  "
$bodyA2  = "In this case the set used to store the result in the layer is empty."
$blankP  = ""
$h2Text  = "2.2.2 The 1st layer"
$bodyB1  = "The 1st layer is populated "
$bodyB2  = "by iterating over all of the cities and as we do so we retrieve each entry in the 0th layer.  If the entry already includes the city being considered, then that entry is skipped.  Otherwise, a new subproblem is added to the 1st layer.  This subproblem consists of the path from the city under consideration to the previous entry" + $apos + "s city, as well as the added cost of the new connection.  This new subproblem is accessed by the tuple of the considered city and the union of the previous entry" + $apos + "s set and the previous entry" + $apos + "s city."

$bigText = $h1Text + $cr + $bodyA1 + $bodyA2 + $cr + $blankP + $cr + $h2Text + $cr + $bodyB1 + $bodyB2 + $cr
$anchorIns.InsertBefore($bigText)

for ($i = 19; $i -le 27; $i++) {
    Write-Host ($i.ToString() + ": [" + $d.Paragraphs.Item($i).Range.Text + "]")
}

# ---------------------------------------------------------------------------
# Step 3: fix paragraph styles for the two new headings.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(20).Style = "IOPH2"
$d.Paragraphs.Item(23).Style = "IOPH2"

Write-Host "Styles fixed"
